$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (Coin name) is untouched text -- used as a style donor to strip the
# temporary Text number-format back off of D-column cells after the literal
# (dotted/leading-zero) numeric-looking strings are typed in, so the cells keep
# their original General-format style index.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.310.94"
$ws.Range("D2").Style = $ws.Range("B2").Style
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.869.63"
$ws.Range("D3").Style = $ws.Range("B3").Style
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.01"
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4706"
$ws.Range("D7").Style = $ws.Range("B7").Style
$ws.Range("E7").Value = "  +0.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2888"
$ws.Range("D8").Style = $ws.Range("B8").Style
$ws.Range("E8").Value = "  +1.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06620"
$ws.Range("D9").Style = $ws.Range("B9").Style
$ws.Range("E9").Value = "  +1.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.70"
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08040"
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.41"
$ws.Range("D12").Style = $ws.Range("B12").Style
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.872.69"
$ws.Range("D13").Style = $ws.Range("B13").Style
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.140"
$ws.Range("D14").Style = $ws.Range("B14").Style
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6867"
$ws.Range("D15").Style = $ws.Range("B15").Style
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "272.31"
$ws.Range("D16").Style = $ws.Range("B16").Style
$ws.Range("E16").Value = "  -2.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.302.99"
$ws.Range("D17").Style = $ws.Range("B17").Style
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.13"
$ws.Range("D18").Style = $ws.Range("B18").Style
$ws.Range("E18").Value = "  +3.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007752"
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").Value = "  +6.10%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.118.65"
$ws.Range("D21").Style = $ws.Range("B21").Style
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.315"
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.220"
$ws.Range("D24").Style = $ws.Range("B24").Style
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.403"
$ws.Range("D25").Style = $ws.Range("B25").Style
$ws.Range("E25").Value = "  +2.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.46"
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.96"
$ws.Range("D27").Style = $ws.Range("B27").Style
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.956"
$ws.Range("D28").Style = $ws.Range("B28").Style
$ws.Range("E28").Value = "  +1.89%  "
$ws.Range("E29").Value = "  -0.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09904"
$ws.Range("D30").Style = $ws.Range("B30").Style
$ws.Range("E30").Value = "  +2.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.376"
$ws.Range("D31").Style = $ws.Range("B31").Style
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.079"
$ws.Range("D33").Style = $ws.Range("B33").Style
$ws.Range("E33").Value = "  +1.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.133"
$ws.Range("D35").Style = $ws.Range("B35").Style
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7021"
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.703"
$ws.Range("D37").Style = $ws.Range("B37").Style
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01888"
$ws.Range("D38").Style = $ws.Range("B38").Style
$ws.Range("E38").Value = "  +1.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.658"
$ws.Range("D39").Style = $ws.Range("B39").Style
$ws.Range("E39").Value = "  +2.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.314"
$ws.Range("D40").Style = $ws.Range("B40").Style
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.76"
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").Value = "  -3.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.958"
$ws.Range("D42").Style = $ws.Range("B42").Style
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8444"
$ws.Range("D43").Style = $ws.Range("B43").Style
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4171"
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.40"
$ws.Range("D46").Style = $ws.Range("B46").Style
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.278"
$ws.Range("D47").Style = $ws.Range("B47").Style
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.088"
$ws.Range("D48").Style = $ws.Range("B48").Style
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "928.10"
$ws.Range("D49").Style = $ws.Range("B49").Style
$ws.Range("E49").Value = "  -4.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.49"
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05681"
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").Value = "  +0.57%  "
